$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(17, 30, 0.01, 0.003, "Regular", "<function relu at 0x1104e69d8>", 0.8981999754905701, 0.3424000144004822, 0.04809999838471413, 0.3424164950847626, 3.330149173736572, 0.3424000144004822, "logs/results_121.log", "weights/model_121.ckpt", "tb/121/non_robust", "(1.4056362, 3.0716696, 6.664684, 11.19474, 13.153981, 10.241105, 6.051113)", "(24.022184, 14.81685, 21.985281, 18.528912, 16.116236, 15.134631, 11.263296, 22.44767)"),
    @(18, 30, 0.02, 0.0003, "Regular", "<function relu at 0x10df1b9d8>", 0.9559000134468079, 0.2694000005722046, 0.05429999902844429, 0.1658958792686462, 4.698654651641846, 0.2694000005722046, "logs/results_123.log", "weights/model_123.ckpt", "tb/123/non_robust", "(0.3431952, 0.41664022, 0.5683927, 0.91492766, 1.3893795, 2.3776739, 3.8029735)", "(6.3056436, 7.597117, 7.726599, 8.445029, 7.708117, 7.5270963, 7.9869733, 9.796146)"),
    @(19, 30, 0.02, 0.0003, "Regular", "<function relu at 0x119b5f9d8>", 0.8652999997138977, 0.1185000017285347, 0.005900000222027302, 0.4752624034881592, 3.402863264083862, 0.1185000017285347, "logs/results_126.log", "weights/model_126.ckpt", "tb/126/non_robust", "(0.90816385, 0.7060737, 0.9453165, 1.4311132, 1.8074349, 2.626249, 3.3506863)", "(15.063638, 2.4488645, 2.8450708, 3.4224455, 2.870352, 2.775072, 2.7803397, 3.6928034)")
)

$startRow = 19
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($r, $col).Value = $data[$c]
    }
}
